$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format for column D so numeric-looking price strings
# (e.g. "0.999") are stored as text, matching the original inlineStr cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.405.19'
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").Value = '3.488.15'
$ws.Range("E3").Value = '  -2.36%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '609.52'
$ws.Range("E5").Value = '  +4.68%  '
$ws.Range("D6").Value = '185.89'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '0.212'
$ws.Range("E9").Value = '  -5.49%  '
$ws.Range("D10").Value = '0.645'
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("D11").Value = '52.84'
$ws.Range("E11").Value = '  -3.34%  '
$ws.Range("D12").Value = '0.0000306'
$ws.Range("E12").Value = '  -4.19%  '
$ws.Range("D13").Value = '9.48'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").Value = '4.033.87'
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("D15").Value = '599.84'
$ws.Range("E15").Value = '  +4.34%  '
$ws.Range("D16").Value = '69.476.77'
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").Value = '18.81'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").Value = '12.56'
$ws.Range("E18").Value = '  -2.03%  '
$ws.Range("D19").Value = '3.499.61'
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("D22").Value = '17.18'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("D23").Value = '105.42'
$ws.Range("E23").Value = '  +12.13%  '
$ws.Range("E24").Value = '  +3.91%  '
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("E26").Value = '  +2.92%  '
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("D28").Value = '9.66'
$ws.Range("E28").Value = '  +4.64%  '
$ws.Range("D29").Value = '33.35'
$ws.Range("E29").Value = '  +2.78%  '
$ws.Range("D30").Value = '6.93'
$ws.Range("E30").Value = '  -3.83%  '
$ws.Range("D31").Value = '4.13'
$ws.Range("E31").Value = '  +14.71%  '
$ws.Range("D32").Value = '12.39'
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("D34").Value = '63.29'
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  -6.92%  '
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  +7.51%  '
$ws.Range("D38").Value = '520.11'
$ws.Range("E38").Value = '  -5.25%  '
$ws.Range("E39").Value = '  -4.91%  '
$ws.Range("D40").Value = '3.582.36'
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").Value = '36.71'
$ws.Range("E41").Value = '  -3.83%  '
$ws.Range("E42").Value = '  -3.27%  '
$ws.Range("E43").Value = '  -0.95%  '
$ws.Range("D44").Value = '0.0461'
$ws.Range("E44").Value = '  -1.57%  '
$ws.Range("D45").Value = '2.94'
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("E46").Value = '  +2.80%  '
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  -4.50%  '
$ws.Range("D48").Value = '8.76'
$ws.Range("E48").Value = '  -6.31%  '
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '130.68'
$ws.Range("E50").Value = '  -2.72%  '
$ws.Range("B51").Value = 'OceanProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range("D51").Value = '1.35'
$ws.Range("E51").Value = '  -9.56%  '
